# Update the "Best Funds" sheet with the new fund performance data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Best Funds")

$data = @(
    @("2022-08", "Gohen", 13.28566009764734, 585569.7282415634, 8630204.692416647, -264954.5912859503),
    @("2022-09", "Trustmind", 2.88726404923217, 10620015.04946109, 43758645.49386902, -2475842.789794324),
    @("2022-10", "Wallington", 3.62743153784914, 41969174.17637323, 196905682.1293281, -2696201.928094937),
    @("2022-11", "Leeder", 7.028914070404999, 2292900.69230508, 26801108.0777934, -8391605.447303778),
    @("2022-12", "Gohen", 77.7271359243389, 205432.9423791791, 16104397.4604946, 68749.71752790874),
    @("2023-01", "Magnum", 7.368744008398862, 7696183.77063203, 59924748.15551724, 4482643.662496128),
    @("2023-02", "Applebead", 12.86993986056393, 9386135.624031307, 129673836.4555226, 511300.1728883917),
    @("2023-03", "Catalysm", 11.49391751539593, 193275.2141059854, 2618183.302886812, -203418.7200761429),
    @("2023-04", "Gohen", 829.6124356178589, 4968.622816473871, 3436686.273461192, 690313.6257966359),
    @("2023-05", "Gohen", 19.81385261671864, 227211.4984409648, 5217218.028098187, -488071.3867241474),
    @("2023-06", "Gohen", 19.69453756641664, 536075.8121401727, 10912906.70472738, 180934.3280547339),
    @("2023-07", "Applebead", 66.01470618098682, 10094131.24441911, 675074648.1844716, 1380591.312594123),
    @("2023-08", "Applebead", 9.516932644869431, 5513343.953973653, 60651817.36173481, -2668350.349795799)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $row++
}
